$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the "Uhrzeit" (time) column B entirely - this shifts every
# following column (Schul/Uni, Adresse, Stadt, Bundesland, PLZ, Tische,
# Teilnehmer) one position to the left.
$ws.Columns.Item(2).Delete() | Out-Null

# The (now) B2 cell holds the "Schul/Uni" entry - update its value from the
# old abbreviation "HKA" to the full name, and drop the inherited border
# formatting so it matches a freshly-typed cell.
$ws.Range("B2").Value = "Hochschule Karlsruhe"
$ws.Range("B2").ClearFormats() | Out-Null

# Reflect where the user was working afterwards.
$ws.Range("B9").Select() | Out-Null
